# The document contains a Word "complex field" ( { m:'Mona_Lisa.jpg'.asImage()
# .setConserveRatio(false).setWidth(100) } ) stored as fldChar begin/instrText/
# fldChar end runs. The parser was switched to TokenIteratorFieldRewriterSplit,
# which expects the field to be represented as literal template text
# ( <w:t>{ ... }</w:t> runs ) instead of a real Word field. Rewrite that
# paragraph accordingly, keeping every run's formatting (rPr) and the
# bookmarkStart/bookmarkEnd untouched.

$d = $word.ActiveDocument

# Locate the paragraph that holds the field (robust to paragraph indices).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find the paragraph containing the field."
}

# Replacement content for the paragraph: the fldChar begin + leading
# " " instrText collapse into a literal "{" run, each instrText run becomes
# an equivalent w:t run (formatting preserved), and the trailing " "
# instrText + fldChar end collapse into a literal "}" run.
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>'</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>Mona_Lisa</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>.jpg</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>'.asImage()</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>.setConserveRatio(false</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>.setWidth(100)</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetPara.Range.InsertXML($xml)
